$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.407.67'
$ws.Range('E2').Value = '  +3.14%  '
$ws.Range('D3').Value = '2.305.37'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.14'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.30'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.43%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.42%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +8.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.59'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0811'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '51.96'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.49%  '
$ws.Range('D15').Value = '2.663.27'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.11'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.86%  '
$ws.Range('D17').Value = '2.309.85'
$ws.Range('E17').Value = '  +1.94%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.810'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.88%  '
$ws.Range('D19').Value = '43.298.62'
$ws.Range('E19').Value = '  +3.33%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').Value = '0.0₃0926'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('E22').Value = '  +3.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.13'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '243.06'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.92%  '
$ws.Range('E25').Value = '  +3.44%  '
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.74'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.27%  '
$ws.Range('E29').Value = '  +8.07%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '37.03'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.66'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '168.51'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.29'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.11'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('E36').Value = '  +6.09%  '
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.87'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.55%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.106'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.47%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.47'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +8.31%  '
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0294'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.30%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.982.66'
$ws.Range('E44').Value = '  +1.60%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '19.00'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.01'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.19%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '56.02'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.21%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.93'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.59'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +9.08%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.532.10'
$ws.Range('E51').Value = '  +2.02%  '
